# [PHOENIX-6082] completed Forward/Close Grievance
#
# Update the trade commencement date on the "tradeDetails" sheet and make
# that sheet the active one (mirroring the user navigating to it and
# editing the date field), moving the selection away from the
# "searchTradeDeatils" sheet that used to be active.

$wb = $excel.ActiveWorkbook

$wsTradeDetails = $wb.Worksheets.Item("tradeDetails")

# Update the trade commencement date value (column H, row 2).
$wsTradeDetails.Range("H2").Value = "31/03/2017"

# Activate the tradeDetails sheet and move the selection to E9, matching
# where the user left off after completing the edit.
$wsTradeDetails.Activate() | Out-Null
$wsTradeDetails.Range("E9").Select() | Out-Null
